$wb = $excel.ActiveWorkbook

# --- Sheet: Costs and Revenues ---
$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Range("B2").Value = 77388.66797673712
$ws.Range("D2").Value = 9992.97670278544
$ws.Range("E2").Value = 9770
$ws.Range("F2").Value = 51915.03718374025

# --- Sheet: Capacities ---
$ws = $wb.Worksheets.Item("Capacities")
$ws.Range("C3").Value = 68

# --- Sheet: PV Dispatch ---
$ws = $wb.Worksheets.Item("PV Dispatch")
$ws.Range("G2").Value = 13.6
$ws.Range("H2").Value = 27.2
$ws.Range("I2").Value = 34
$ws.Range("J2").Value = 40.8
$ws.Range("K2").Value = 47.6
$ws.Range("L2").Value = 54.4
$ws.Range("M2").Value = 61.2
$ws.Range("N2").Value = 68
$ws.Range("O2").Value = 61.2
$ws.Range("P2").Value = 54.4
$ws.Range("Q2").Value = 47.6
$ws.Range("R2").Value = 34
$ws.Range("S2").Value = 20.4
$ws.Range("T2").Value = 13.6
$ws.Range("I3").Value = 27.2
$ws.Range("J3").Value = 40.8
$ws.Range("K3").Value = 54.4
$ws.Range("L3").Value = 61.2
$ws.Range("M3").Value = 68
$ws.Range("N3").Value = 54.4
$ws.Range("O3").Value = 47.6
$ws.Range("P3").Value = 34
$ws.Range("Q3").Value = 34
$ws.Range("R3").Value = 20.4
$ws.Range("S3").Value = 13.6
$ws.Range("K4").Value = 27.2
$ws.Range("L4").Value = 47.6
$ws.Range("M4").Value = 51.18312417100189
$ws.Range("N4").Value = 54.4
$ws.Range("O4").Value = 47.6
$ws.Range("P4").Value = 27.2
$ws.Range("Q4").Value = 13.6

# --- Sheet: Battery Input ---
$ws = $wb.Worksheets.Item("Battery Input")
$ws.Range("G2").Value = 64.3
$ws.Range("H2").Value = 14.2
$ws.Range("I2").Value = 2.8
$ws.Range("J2").Value = 60.3
$ws.Range("K2").Value = 21.6
$ws.Range("L2").Value = 33.6
$ws.Range("M2").Value = 37.8
$ws.Range("N2").Value = 42
$ws.Range("O2").Value = 30
$ws.Range("P2").Value = 25.8
$ws.Range("Q2").Value = 73.53427201306106
$ws.Range("R2").Value = 0.2
$ws.Range("S2").Value = 32.4
$ws.Range("T2").Value = 45.6
$ws.Range("I3").Value = 27.2
$ws.Range("J3").Value = 40.8
$ws.Range("K3").Value = 54.4
$ws.Range("L3").Value = 61.2
$ws.Range("M3").Value = 44.6
$ws.Range("N3").Value = 28.4
$ws.Range("O3").Value = 47.6
$ws.Range("P3").Value = 5.4
$ws.Range("Q3").Value = 8.230792776247645
$ws.Range("R3").Value = 20.4
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 27.2
$ws.Range("L4").Value = 47.6
$ws.Range("M4").Value = 27.78312417100189
$ws.Range("N4").Value = 54.4
$ws.Range("O4").Value = 47.6
$ws.Range("P4").Value = 27.2
$ws.Range("Q4").Value = 13.6
$ws.Range("R4").Value = 0

# --- Sheet: State of Charge ---
$ws = $wb.Worksheets.Item("State of Charge")
$ws.Range("G2").Value = 183.657
$ws.Range("H2").Value = 197.715
$ws.Range("I2").Value = 200.487
$ws.Range("J2").Value = 260.184
$ws.Range("K2").Value = 281.568
$ws.Range("L2").Value = 314.832
$ws.Range("M2").Value = 352.254
$ws.Range("N2").Value = 393.834
$ws.Range("O2").Value = 423.534
$ws.Range("P2").Value = 449.076
$ws.Range("Q2").Value = 521.8749292929305
$ws.Range("R2").Value = 522.0729292929304
$ws.Range("S2").Value = 554.1489292929305
$ws.Range("I3").Value = 146.928
$ws.Range("J3").Value = 187.32
$ws.Range("K3").Value = 241.176
$ws.Range("L3").Value = 301.764
$ws.Range("M3").Value = 345.918
$ws.Range("N3").Value = 374.034
$ws.Range("O3").Value = 421.158
$ws.Range("P3").Value = 426.504
$ws.Range("Q3").Value = 434.6524848484852
$ws.Range("J4").Value = 120
$ws.Range("K4").Value = 146.928
$ws.Range("L4").Value = 194.052
$ws.Range("M4").Value = 221.5572929292919
$ws.Range("N4").Value = 275.4132929292919
$ws.Range("O4").Value = 322.5372929292918
$ws.Range("P4").Value = 349.4652929292918
$ws.Range("Q4").Value = 362.9292929292918

# --- Sheet: Feed in from Type 2 ---
$ws = $wb.Worksheets.Item("Feed in from Type 2")
$ws.Range("J2").Value = 0
$ws.Range("Q2").Value = 51.93427201306106
$ws.Range("T2").Value = 34
$ws.Range("N3").Value = 0
$ws.Range("J4").Value = 0

# --- Sheet: Feed in from Type 3 ---
$ws = $wb.Worksheets.Item("Feed in from Type 3")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0

# --- Sheet: Feed in from Type 4 ---
$ws = $wb.Worksheets.Item("Feed in from Type 4")
$ws.Range("J2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("Q3").Value = 0.2307927762476449
$ws.Range("S3").Value = 9.6
$ws.Range("J4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
